$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "First line from person1"
$ws.Range("C6").Select()
